$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from Sheet1 to Sheet2
$ws.Name = "Sheet2"

# Update the label text in A2 (shared string "No control" -> "Holding control")
$ws.Range("A2").Value = "Holding control"

# Update the numeric values in row 2
$ws.Range("B2").Value = 2450.824609460746
$ws.Range("C2").Value = 12623.373108514
$ws.Range("D2").Value = 9.034912022317402
$ws.Range("E2").Value = 15083.23262999707
